$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (logistic_embeddings)
$ws.Range("C5").Value = 0.396
$ws.Range("D5").Value = 0.498
$ws.Range("E5").Value = 0.534
$ws.Range("F5").Value = 0.588
$ws.Range("G5").Value = 0.57
$ws.Range("H5").Value = 0.601

# Row 7 (classical-best-embeddings -> classical-best-embed)
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.396
$ws.Range("E7").Value = 0.534
$ws.Range("F7").Value = 0.588
$ws.Range("H7").Value = 0.601

# Row 8 (BERT-base)
$ws.Range("C8").Value = 0.403
$ws.Range("D8").Value = 0.5669999999999999
$ws.Range("E8").Value = 0.602
$ws.Range("F8").Value = 0.638
$ws.Range("G8").Value = 0.634
$ws.Range("H8").Value = 0.666

# Row 9 (BERT-base-nli)
$ws.Range("B9").Value = 0.46
$ws.Range("C9").Value = 0.495
$ws.Range("D9").Value = 0.574
$ws.Range("E9").Value = 0.6
$ws.Range("G9").Value = 0.577
$ws.Range("H9").Value = 0.609
